$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New matchup rows (fall 24 week 13 inputs) appended after existing data (row 1946)
$data = @(
    @(6,18,5,2),
    @(4,14,3,6),
    @(4,4,6,16),
    @(4,17,6,3),
    @(5,17,4,3),
    @(5,7,7,13),
    @(5,14,4,6),
    @(3,7,4,13),
    @(4,4,2,16),
    @(4,3,3,17),
    @(4,15,5,5),
    @(7,8,6,12),
    @(4,14,3,6),
    @(4,16,5,4),
    @(5,19,4,1),
    @(5,4,4,16),
    @(5,4,6,16),
    @(7,3,5,17),
    @(2,14,4,6),
    @(4,7,5,13),
    @(4,8,3,12),
    @(6,14,9,6),
    @(7,7,5,13),
    @(4,14,5,6),
    @(5,13,4,7),
    @(3,13,5,7),
    @(6,14,5,6),
    @(5,17,3,3),
    @(6,14,7,6),
    @(3,8,4,12),
    @(6,8,3,12),
    @(3,6,4,14),
    @(4,5,3,15),
    @(5,16,4,4),
    @(4,6,5,14),
    @(4,15,2,5),
    @(4,6,3,14),
    @(6,4,4,16),
    @(4,12,7,8),
    @(4,13,3,7)
)

$startRow = 1947
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update view/selection to match post-edit state
$ws.Application.ActiveWindow.ScrollRow = 1970
$newSelection = "A" + ($endRow + 1)
$ws.Range($newSelection).Select()
